# "Generate Report for Handback"
#
# Overview sheet: the Status text shown for the sample file flips from
# "Ready for handoff" to "Handed back: in sync with en-US" (columns E + F).
# The same shared text also backs the per-language "Status" column (C) on
# the zh-cn / de-de sheets, so setting it once per occurrence lets the
# engine's shared-string dedup collapse them back down to a single string,
# just like the source diff shows.
#
# zh-cn / de-de sheets: the handback has now produced a "Latest Target
# File" (I), a "Latest Handback File" (J) and a "Latest Handback DateTime"
# (K) for each locale, plus a new hyperlink on I2 pointing at the handed-
# back markdown file (mirroring the existing hyperlink on A2).
#
# A few columns are also widened to better fit the new long filenames.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$mdFile = "91171a96-5f7d-4641-8576-44b0058bbe25.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1684b09a54dcde7f52c067f3dd5e276c8aa47a75/e2e/91171a96-5f7d-4641-8576-44b0058bbe25.md"

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Column widths grew to fit the longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666664
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666664

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus

$wsZhCn.Range("I2").Value = $mdFile
$wsZhCn.Range("J2").Value = "91171a96-5f7d-4641-8576-44b0058bbe25.9800bf1d99ddd4eb7800e02a0ef556eabf3e0ef6.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-29 02:58:34"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFile)
$wsZhCn.Range("I2").Font.Underline = 2
$wsZhCn.Range("I2").Font.Color = 15570276

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666664
$wsZhCn.Columns.Item(9).ColumnWidth = 39.16666666666667
$wsZhCn.Columns.Item(10).ColumnWidth = 39.16666666666667

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus

$wsDeDe.Range("I2").Value = $mdFile
$wsDeDe.Range("J2").Value = "91171a96-5f7d-4641-8576-44b0058bbe25.9800bf1d99ddd4eb7800e02a0ef556eabf3e0ef6.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-29 02:58:41"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFile)
$wsDeDe.Range("I2").Font.Underline = 2
$wsDeDe.Range("I2").Font.Color = 15570276

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666664
$wsDeDe.Columns.Item(9).ColumnWidth = 39.16666666666667
$wsDeDe.Columns.Item(10).ColumnWidth = 39.16666666666667
